$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns BM:BO (school2_4yr, school2_ad, school2_cc).
# This shifts BP:BV left by 3 into BM:BS.
$ws.Range("BM:BO").Delete()

# Rename BL header and recode its data to 0.
$ws.Range("BL1").Value = "school2_type"
$ws.Range("BL2:BL61").Value = 0
